# Apply hybrid bold + color highlighting to quantitative impact metrics
# across the resume bullet points, matching the target diff.
#
# Color used: RGB(0x2C, 0x3E, 0x50) -> packed BGR integer for Word's
# Font.Color property = 0x2C + (0x3E * 256) + (0x50 * 65536) = 5258796

$d = $word.ActiveDocument
$highlightColor = 5258796

function Highlight-InParagraph {
    param(
        [int]$paraIndex,
        [string[]]$terms
    )

    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range.Duplicate

    foreach ($term in $terms) {
        $found = $rng.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $highlightColor
            # Continue searching in the remainder of the paragraph after this match
            $rng.Start = $rng.End
            $rng.End = $para.Range.End
        }
    }
}

# NOTE: $d.Paragraphs(N) is 1-based.

# Para 10: "Discovered systematic race coding errors ... from 23% to 64%"
Highlight-InParagraph 10 @("23%", "64%")

# Para 12: "Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
Highlight-InParagraph 12 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%", "71%", "87%")

# Para 13: "Trigonometric algorithm ... by 73.5%, saving ... organizations $4.7M ..."
Highlight-InParagraph 13 @("73.5%", "$4.7M")

# Para 14: "Built real-time FEC analysis systems ... valued over $2 trillion"
Highlight-InParagraph 14 @("$2")

# Para 24: "Modernized legacy ETL processes ... reducing processing time by 57%"
Highlight-InParagraph 24 @("57%")

# Para 50: "Revenue generation: Delivered $4.9M additional revenue through optimization"
Highlight-InParagraph 50 @("$4.9M")

# Para 51: "23% conversion rate improvement"
Highlight-InParagraph 51 @("23%")

# Para 53: "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Highlight-InParagraph 53 @("12,847")
